# Calc-Results.xlsx commit: "data, gnuplot euler, update each code for exact solution"
#
# - rename Sheet1 -> IVP2
# - move selection from I19 to I9
# - replace the computed-result column (G3:G20) with the new values that
#   came out of re-running the Euler / Heun / RK4 code against the exact
#   solution (x - x*x)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "IVP2"

# New result values for column G (rows 3-20)
$ws.Range("G3").Value  = 0.064161886868
$ws.Range("G4").Value  = 0.0063939021304
$ws.Range("G5").Value  = 0.0006386369113
$ws.Range("G6").Value  = 0.43526659839
$ws.Range("G7").Value  = 0.019467951076
$ws.Range("G8").Value  = 0.0018493303035
$ws.Range("G9").Value  = 0.009163670038
$ws.Range("G10").Value = 0.00010831145445
$ws.Range("G11").Value = 0.0000011042980272
$ws.Range("G12").Value = 0.064733401606
$ws.Range("G13").Value = 0.00064334177177
$ws.Range("G14").Value = 0.000006055373694
$ws.Range("G15").Value = 0.000063326523133
$ws.Range("G16").Value = 0.0000000076208883648
$ws.Range("G17").Value = 0.00000000000077637896112
$ws.Range("G18").Value = 0.00034273915799
$ws.Range("G19").Value = 0.0000002696915673
$ws.Range("G20").Value = 0.000000000026141755427

# Keep the scientific-notation display format on the results column
$ws.Range("G3:G20").NumberFormat = "0.0000000000E+00"

# Move the active selection to I9 (was I19)
$ws.Range("I9").Select()
